$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.170.05"
$ws.Range("E2").Value = "  +1.79%  "

$ws.Range("D3").Value = "'1.909.51"
$ws.Range("E3").Value = "  +2.04%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'327.88"
$ws.Range("E5").Value = "  +0.44%  "

$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").Value = "'0.4647"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").Value = "'0.3926"
$ws.Range("E8").Value = "  +1.05%  "

$ws.Range("D9").Value = "'46.98"
$ws.Range("E9").Value = "  +0.87%  "

$ws.Range("D10").Value = "'0.07963"
$ws.Range("E10").Value = "  +1.11%  "

$ws.Range("E11").Value = "  +2.97%  "

$ws.Range("D12").Value = "'22.28"
$ws.Range("E12").Value = "  +1.40%  "

$ws.Range("D13").Value = "'1.924.97"
$ws.Range("E13").Value = "  +0.58%  "

$ws.Range("D14").Value = "'7.133"
$ws.Range("E14").Value = "  +2.03%  "

$ws.Range("D15").Value = "'5.793"
$ws.Range("E15").Value = "  +1.54%  "

$ws.Range("D16").Value = "'0.06983"
$ws.Range("E16").Value = "  -0.19%  "

$ws.Range("D17").Value = "'88.64"
$ws.Range("E17").Value = "  +0.60%  "

$ws.Range("D18").Value = "'1.005"
$ws.Range("E18").Value = "  +0.01%  "

$ws.Range("D19").Value = "'0.00001011"
$ws.Range("E19").Value = "  +0.61%  "

$ws.Range("D20").Value = "'17.27"
$ws.Range("E20").Value = "  +2.58%  "

$ws.Range("D21").Value = "'1.003"
$ws.Range("E21").Value = "  -0.08%  "

$ws.Range("D22").Value = "'29.190.33"
$ws.Range("E22").Value = "  +1.88%  "

$ws.Range("D23").Value = "'5.373"
$ws.Range("E23").Value = "  +1.59%  "

$ws.Range("D24").Value = "'11.08"
$ws.Range("E24").Value = "  +0.64%  "

$ws.Range("D25").Value = "'2.173.48"
$ws.Range("E25").Value = "  +1.77%  "

$ws.Range("D26").Value = "'2.058"
$ws.Range("E26").Value = "  -2.85%  "

$ws.Range("D27").Value = "'156.34"
$ws.Range("E27").Value = "  +2.28%  "

$ws.Range("E28").Value = "  +1.66%  "

$ws.Range("D29").Value = "'5.843"
$ws.Range("E29").Value = "  +1.04%  "

$ws.Range("D30").Value = "'2.003"
$ws.Range("E30").Value = "  +0.74%  "

$ws.Range("D31").Value = "'119.63"
$ws.Range("E31").Value = "  +0.18%  "

$ws.Range("D32").Value = "'0.09399"
$ws.Range("E32").Value = "  +0.35%  "

$ws.Range("D33").Value = "'0.9233"
$ws.Range("E33").Value = "  +0.27%  "

$ws.Range("E34").Value = "  +1.87%  "

$ws.Range("D35").Value = "'1.344"
$ws.Range("E35").Value = "  +0.40%  "

$ws.Range("D36").Value = "'3.277"
$ws.Range("E36").Value = "  -1.76%  "

$ws.Range("D37").Value = "'0.05842"
$ws.Range("E37").Value = "  +0.75%  "

$ws.Range("D38").Value = "'1.160"
$ws.Range("E38").Value = "  +0.96%  "

$ws.Range("D39").Value = "'8.020"

$ws.Range("D40").Value = "'0.02097"
$ws.Range("E40").Value = "  -0.16%  "

$ws.Range("D41").Value = "'0.5755"
$ws.Range("E41").Value = "  +2.24%  "

$ws.Range("E42").Value = "  +1.19%  "

$ws.Range("D43").Value = "'10.000"
$ws.Range("E43").Value = "  +2.35%  "

$ws.Range("D44").Value = "'12.00"
$ws.Range("E44").Value = "  +2.80%  "

$ws.Range("D45").Value = "'0.5426"
$ws.Range("E45").Value = "  +2.08%  "

$ws.Range("D46").Value = "'2.219"
$ws.Range("E46").Value = "  +5.90%  "

$ws.Range("D47").Value = "'0.07098"
$ws.Range("E47").Value = "  -1.57%  "

$ws.Range("D48").Value = "'1.885"
$ws.Range("E48").Value = "  +3.18%  "

$ws.Range("D49").Value = "'2.590"
$ws.Range("E49").Value = "  +6.89%  "

$ws.Range("E50").Value = "  -0.82%  "

$ws.Range("D51").Value = "'1.075"
$ws.Range("E51").Value = "  -6.27%  "

